function Set-ShapeText($shape, $text) {
    # Preserve the shape's original size - these textboxes use spAutoFit,
    # and PowerPoint will otherwise resize the box to fit the new text.
    $origH = $shape.Height
    $origW = $shape.Width
    $shape.TextFrame.TextRange.Text = $text
    $shape.Height = $origH
    $shape.Width = $origW
}

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# Slide 2 (index 2) - "Cost Comparison Summary"
# ---------------------------------------------------------------------------
$s2 = $p.Slides.Item(2)

# Subtitle paragraph (TextBox 4)
Set-ShapeText $s2.Shapes.Item(4) "Replacing the Waracle Affordability Calculator. Current Waracle support contract (£3,775/month) ends February 2025. In-house go-live targeted for end of February, enabling seamless transition with no renewal required."

# Table row label "Annual Support & Hosting" -> add "(replaces Waracle)"
$tbl2 = $s2.Shapes.Item(5).Table
$tbl2.Cell(3, 1).Shape.TextFrame.TextRange.Text = "Annual Support & Hosting (replaces Waracle)"

# Bottom callout text (TextBox 8)
Set-ShapeText $s2.Shapes.Item(8) "Saves £3,775/month (£45,300/year) from Waracle contract not renewed, plus avoids £120,000/year Podium costs"

# ---------------------------------------------------------------------------
# Slide 4 (index 4) - "External Costs Avoided"
# ---------------------------------------------------------------------------
$s4 = $p.Slides.Item(4)

# Subtitle paragraph (TextBox 4)
Set-ShapeText $s4.Shapes.Item(4) "Current Waracle contract: £3,775/month (£45,300/year). Contract ends February 2025. Podium would have replaced Waracle at higher cost."

# Table rows (Table 5)
$tbl4 = $s4.Shapes.Item(5).Table
$tbl4.Cell(3, 1).Shape.TextFrame.TextRange.Text = "Monthly support (replaces Waracle)"
$tbl4.Cell(3, 4).Shape.TextFrame.TextRange.Text = "vs current £3,775/month"
$tbl4.Cell(4, 1).Shape.TextFrame.TextRange.Text = "Annual support cost"
$tbl4.Cell(4, 4).Shape.TextFrame.TextRange.Text = "vs current £45,300/year"

# ---------------------------------------------------------------------------
# Slide 5 (index 5) - "Cumulative Savings Over Time"
# ---------------------------------------------------------------------------
$s5 = $p.Slides.Item(5)

# Subtitle paragraph (TextBox 4)
Set-ShapeText $s5.Shapes.Item(4) "In-house go-live end of February 2025 enables Waracle contract to end without renewal. Comparison vs Podium alternative (realistic estimates)."

# ---------------------------------------------------------------------------
# Slide 6 (index 6) - "Summary and Benefits"
# ---------------------------------------------------------------------------
$s6 = $p.Slides.Item(6)

# Table rows (Table 4)
$tbl6 = $s6.Shapes.Item(4).Table
$tbl6.Cell(5, 1).Shape.TextFrame.TextRange.Text = "Waracle contract saved (from March)"
$tbl6.Cell(5, 2).Shape.TextFrame.TextRange.Text = "£3,775/month (£45,300/year)"

# Additional benefit bullet textboxes - cascading shift, new bullet inserted at top
Set-ShapeText $s6.Shapes.Item(6) "• Waracle contract ends February 2025 - no renewal required"
Set-ShapeText $s6.Shapes.Item(7) "• No ongoing third-party hosting or support costs"
Set-ShapeText $s6.Shapes.Item(8) "• No external dependency for future changes or FSMA compliance"
Set-ShapeText $s6.Shapes.Item(9) "• Internal capability and knowledge retention"
Set-ShapeText $s6.Shapes.Item(10) "• Full control over codebase and future roadmap"
